$d = $word.ActiveDocument

# Direct Range.Text assignment is used (instead of Find.Execute replacement)
# so straight apostrophes and the xml:space="preserve" attribute on runs
# that need it are preserved verbatim, and only the targeted paragraph's
# run text changes while its run formatting (rPr) is kept intact.

$d.Paragraphs.Item(1).Range.Text  = "Le problème des fourmis - sous-titres :"
$d.Paragraphs.Item(7).Range.Text  = "[Musique]"
$d.Paragraphs.Item(11).Range.Text = "Les énigmes que je vous"
$d.Paragraphs.Item(15).Range.Text = "poserai sont des versions élémentaires"
$d.Paragraphs.Item(19).Range.Text = "d'une énigme plus compliquée"
